$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking"): Right column 5 -> 4, Wrong column -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right column 90 -> 72, Wrong column -3 -> -6
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -6

# E12 summary text reflects the corrected totals
$ws.Range("E12").Value = "66 / 112"
